$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove all existing hyperlinks on this sheet so stale relationship targets
# (pointing at old URLs) don't linger under reused cell refs.
$ws.Cells.Hyperlinks.Delete()

# --- Row 1 header stays as-is ---

# --- Rewrite data rows 2-15 ---
# Row 2
$ws.Cells.Item(2, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(2, 2).Value = 'ChatGPTを用いた当事業部内チャットツールのシステム開発'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5431738'
$ws.Cells.Item(2, 7).Value = 463
$ws.Cells.Item(2, 8).Value = '🔥GPT,ChatGPT ◆ツール,開発'
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5431738', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431738') | Out-Null
$ws.Cells.Item(2, 6).Style = "Hyperlink"

# Row 3
$ws.Cells.Item(3, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(3, 2).Value = '【急募】AIシステム構築!FirebaseとOpenAI活用の専門家募集'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5431299'
$ws.Cells.Item(3, 7).Value = 325
$ws.Cells.Item(3, 8).Value = '🔥AI,Ai'
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5431299', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431299') | Out-Null
$ws.Cells.Item(3, 6).Style = "Hyperlink"

# Row 4
$ws.Cells.Item(4, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(4, 2).Value = '【急募】Cordova必須!スマホアプリ開発支援メンバー募集'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5431740'
$ws.Cells.Item(4, 7).Value = 175
$ws.Cells.Item(4, 8).Value = '★スマホアプリ ◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5431740', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431740') | Out-Null
$ws.Cells.Item(4, 6).Style = "Hyperlink"

# Row 5
$ws.Cells.Item(5, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(5, 2).Value = 'webアプリの開発'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5431673'
$ws.Cells.Item(5, 7).Value = 100
$ws.Cells.Item(5, 8).Value = '◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5431673', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431673') | Out-Null
$ws.Cells.Item(5, 6).Style = "Hyperlink"

# Row 6
$ws.Cells.Item(6, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(6, 2).Value = '【急募】Webアプリ開発エンジニア募集!フルリモート可'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5411585'
$ws.Cells.Item(6, 7).Value = 93
$ws.Cells.Item(6, 8).Value = '◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5411585', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5411585') | Out-Null
$ws.Cells.Item(6, 6).Style = "Hyperlink"

# Row 7
$ws.Cells.Item(7, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(7, 2).Value = '【急募】知的財産関連システムの開発パートナー募集'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5431547'
$ws.Cells.Item(7, 7).Value = 90
$ws.Cells.Item(7, 8).Value = '◆開発'
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5431547', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431547') | Out-Null
$ws.Cells.Item(7, 6).Style = "Hyperlink"

# Row 8
$ws.Cells.Item(8, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(8, 2).Value = '進行管理およびチームディレクションを担当'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5418064'
$ws.Cells.Item(8, 7).Value = 30
$ws.Cells.Item(8, 8).Value = '◇管理'
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5418064', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5418064') | Out-Null
$ws.Cells.Item(8, 6).Style = "Hyperlink"

# Row 9
$ws.Cells.Item(9, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(9, 2).Value = '【急募】Laravel12でFortifyを使った2段階認証システムの制作'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5431508'
$ws.Cells.Item(9, 7).Value = 33
$ws.Cells.Item(9, 8).ClearContents()
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5431508', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431508') | Out-Null
$ws.Cells.Item(9, 6).Style = "Hyperlink"

# Row 10
$ws.Cells.Item(10, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(10, 2).Value = '〖リモート可〗Delphiエンジニア募集'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5341051'
$ws.Cells.Item(10, 7).Value = 25
$ws.Cells.Item(10, 8).ClearContents()
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5341051', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5341051') | Out-Null
$ws.Cells.Item(10, 6).Style = "Hyperlink"

# Row 11
$ws.Cells.Item(11, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(11, 2).Value = '【フルリモート】SESエンジニア募集|スキルに応じて30〜40万円/月|複数案件あり・継続前提'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5417644'
$ws.Cells.Item(11, 7).Value = 25
$ws.Cells.Item(11, 8).ClearContents()
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5417644', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5417644') | Out-Null
$ws.Cells.Item(11, 6).Style = "Hyperlink"

# Row 12
$ws.Cells.Item(12, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(12, 2).Value = '【高単価×長期案件あり】フリーランスエンジニア募集|リモート可・週3〜OK'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5431322'
$ws.Cells.Item(12, 7).Value = 25
$ws.Cells.Item(12, 8).ClearContents()
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5431322', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431322') | Out-Null
$ws.Cells.Item(12, 6).Style = "Hyperlink"

# Row 13
$ws.Cells.Item(13, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(13, 2).Value = '初回 Hubspot構築者募集'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5431947'
$ws.Cells.Item(13, 7).Value = 18
$ws.Cells.Item(13, 8).ClearContents()
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), 'https://www.lancers.jp/work/detail/5431947', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431947') | Out-Null
$ws.Cells.Item(13, 6).Style = "Hyperlink"

# Row 14
$ws.Cells.Item(14, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(14, 2).Value = 'AWS環境からAWS環境ヘの新規構築'
$ws.Cells.Item(14, 3).Value = 'システム開発'
$ws.Cells.Item(14, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = '期限情報なし'
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5431069'
$ws.Cells.Item(14, 7).Value = 18
$ws.Cells.Item(14, 8).ClearContents()
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), 'https://www.lancers.jp/work/detail/5431069', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431069') | Out-Null
$ws.Cells.Item(14, 6).Style = "Hyperlink"

# Row 15
$ws.Cells.Item(15, 1).Value = '2025-11-11 12:37:00'
$ws.Cells.Item(15, 2).Value = 'EAの作成'
$ws.Cells.Item(15, 3).Value = 'システム開発'
$ws.Cells.Item(15, 4).Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Cells.Item(15, 5).Value = '期限情報なし'
$ws.Cells.Item(15, 6).Value = 'https://www.lancers.jp/work/detail/5431276'
$ws.Cells.Item(15, 7).Value = 10
$ws.Cells.Item(15, 8).ClearContents()
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), 'https://www.lancers.jp/work/detail/5431276', [Type]::Missing, [Type]::Missing, 'https://www.lancers.jp/work/detail/5431276') | Out-Null
$ws.Cells.Item(15, 6).Style = "Hyperlink"

# --- Column H width 14 -> 22 ---
$ws.Columns.Item(8).ColumnWidth = 21.16

Write-Output "applied"
